# Generate Report for Handoff
# - Update the "Status" text from "Handed back: in sync with en-US" to
#   "Ready for handoff" on the Overview sheet (both language columns) and on
#   each language-specific sheet.
# - Bump the associated report timestamps to reflect the new handoff run.
# - Narrow the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime"
#   columns now that the text they hold is shorter.

$wb = $excel.ActiveWorkbook

$oldStatus = "Handed back: in sync with en-US"
$newStatus = "Ready for handoff"

$newGenDate = "2016-09-07 09:48:46"
$newZhHandoffDate = "2016-09-07 09:48:35"

# Narrower report columns now hold shorter "Ready for handoff" text. Excel's
# ColumnWidth setter snaps to its internal pixel grid, so 16.33 is the
# character-width value whose stored column width lands nearest the target.
$newColumnWidth = 16.33

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("G2").Value = $newGenDate

$wsOverview.Range("E1").EntireColumn.ColumnWidth = $newColumnWidth
$wsOverview.Range("F1").EntireColumn.ColumnWidth = $newColumnWidth

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("H2").Value = $newZhHandoffDate

$wsZh.Range("C1").EntireColumn.ColumnWidth = $newColumnWidth

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("H2").Value = $newGenDate

$wsDe.Range("C1").EntireColumn.ColumnWidth = $newColumnWidth
